# chore: update Sheets via scheduled runner
#
# Refreshes the cached Universalis market-price figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -- columns H:N of
# each job table) with the latest pull. Values only; no structural changes.

$wb = $excel.ActiveWorkbook

$updates = @"
Sheet,Cell,Value
ALC,H19,2381.6428
ALC,I19,2579.1
ALC,J19,1888
ALC,K19,2579.1
ALC,L19,1888
ALC,M19,-2404.1
ALC,N19,-2238
ALC,H53,334.42856
ALC,I53,191
ALC,J53,525.6667
ALC,K53,191
ALC,L53,525.6667
ALC,M53,446
ALC,N53,-1799.6667
ALC,H62,1770.8334
ALC,I62,1509.625
ALC,J62,2293.25
ALC,K62,1509.625
ALC,L62,2293.25
ALC,M62,-885.625
ALC,N62,-3541.25
ALC,H65,1770.8334
ALC,I65,1509.625
ALC,J65,2293.25
ALC,K65,7548.125
ALC,L65,11466.25
ALC,M65,-4428.125
ALC,N65,-17706.25
ALC,H106,2199.75
ALC,I106,2199.75
ALC,J106,0
ALC,K106,2199.75
ALC,L106,0
ALC,M106,-1568.75
ALC,H125,12831.5
ALC,I125,2000
ALC,J125,14997.8
ALC,K125,18000
ALC,L125,134980.2
ALC,M125,-15540
ALC,N125,-139900.2
ARM,H102,1393.25
ARM,I102,1392.6364
ARM,J102,1400
ARM,K102,1392.6364
ARM,L102,1400
ARM,M102,229.3635999999999
ARM,N102,-4644
ARM,H122,1370.7916
ARM,I122,1090.4762
ARM,J122,3333
ARM,K122,3271.4286
ARM,L122,9999
ARM,M122,-821.4286000000002
ARM,N122,-14899
BSM,H76,117056.336
BSM,I76,0
BSM,J76,117056.336
BSM,K76,0
BSM,L76,117056.336
BSM,N76,-117686.336
BSM,H79,117056.336
BSM,I79,0
BSM,J79,117056.336
BSM,K79,0
BSM,L79,117056.336
BSM,N79,-119240.336
BSM,H80,635.375
BSM,I80,939
BSM,J80,453.2
BSM,K80,939
BSM,L80,453.2
BSM,M80,59
BSM,N80,-2449.2
BSM,H83,635.375
BSM,I83,939
BSM,J83,453.2
BSM,K83,4695
BSM,L83,2266
BSM,M83,297
BSM,N83,-12250
BSM,H86,1816.1666
BSM,I86,1774.5
BSM,J86,1899.5
BSM,K86,1774.5
BSM,L86,1899.5
BSM,M86,-651.5
BSM,N86,-4145.5
BSM,H89,1816.1666
BSM,I89,1774.5
BSM,J89,1899.5
BSM,K89,8872.5
BSM,L89,9497.5
BSM,M89,-3256.5
BSM,N89,-20729.5
BSM,H105,2254.0833
BSM,I105,2254.9
BSM,J105,2250
BSM,K105,2254.9
BSM,L105,2250
BSM,M105,-507.9000000000001
BSM,N105,-5744
BSM,H107,2668
BSM,I107,2751.75
BSM,J107,1998
BSM,K107,2751.75
BSM,L107,1998
BSM,M107,-831.75
BSM,N107,-5838
CRP,H58,2066.5625
CRP,I58,1915.0714
CRP,J58,3127
CRP,K58,1915.0714
CRP,L58,3127
CRP,M58,-1712.0714
CRP,N58,-3533
CRP,H86,13427.917
CRP,I86,13767.5
CRP,J86,12748.75
CRP,K86,13767.5
CRP,L86,12748.75
CRP,M86,-12644.5
CRP,N86,-14994.75
CRP,H89,13427.917
CRP,I89,13767.5
CRP,J89,12748.75
CRP,K89,68837.5
CRP,L89,63743.75
CRP,M89,-63221.5
CRP,N89,-74975.75
CRP,H99,2708.25
CRP,I99,2708.25
CRP,J99,0
CRP,K99,2708.25
CRP,L99,0
CRP,M99,
CRP,N99,-1210.25
CRP,H122,2388.2222
CRP,I122,1999.25
CRP,J122,5500
CRP,K122,5997.75
CRP,L122,16500
CRP,M122,-3547.75
CRP,N122,-21400
CRP,H126,2708.25
CRP,I126,2708.25
CRP,J126,0
CRP,K126,8124.75
CRP,L126,0
CRP,M126,
CRP,N126,-5654.75
CRP,H132,4454
CRP,I132,4855.143
CRP,J132,3050
CRP,K132,14565.429
CRP,L132,9150
CRP,M132,-12035.429
CRP,N132,-14210
CRP,H136,2066.5625
CRP,I136,1915.0714
CRP,J136,3127
CRP,K136,5745.2142
CRP,L136,9381
CRP,M136,-3195.2142
CRP,N136,-14481
CUL,H7,7
CUL,I7,7
CUL,J7,0
CUL,K7,21
CUL,L7,0
CUL,M7,91
CUL,H32,3247.5
CUL,I32,0
CUL,J32,3247.5
CUL,K32,0
CUL,L32,
CUL,M32,9742.5
CUL,N32,-10308.5
CUL,H33,799.25
CUL,I33,699.3333
CUL,J33,1099
CUL,K33,4195.9998
CUL,L33,6594
CUL,M33,-3912.9998
CUL,N33,-7160
CUL,H34,539
CUL,I34,539
CUL,J34,0
CUL,K34,1617
CUL,L34,0
CUL,M34,
CUL,N34,-1533
CUL,H41,2099
CUL,I41,2099
CUL,J41,0
CUL,K41,6297
CUL,L41,0
CUL,M41,-5959
CUL,H45,0
CUL,I45,0
CUL,J45,0
CUL,K45,0
CUL,L45,
CUL,N45,0
CUL,H92,422.8
CUL,I92,401
CUL,J92,437.33334
CUL,K92,1203
CUL,L92,1312.00002
CUL,M92,45
CUL,N92,-3808.00002
CUL,H97,1212.6666
CUL,I97,1622
CUL,J97,394
CUL,K97,4866
CUL,L97,1182
CUL,M97,-4370
CUL,N97,-2174
CUL,H128,592572.3
CUL,I128,592572.3
CUL,J128,0
CUL,K128,1777716.9
CUL,L128,0
CUL,M128,-1772736.9
CUL,H134,2659.8
CUL,I134,2659.8
CUL,J134,0
CUL,K134,7979.400000000001
CUL,L134,0
CUL,M134,-2909.400000000001
GSM,H70,5185.6665
GSM,I70,5185.6665
GSM,J70,0
GSM,K70,5185.6665
GSM,L70,0
GSM,M70,-4915.6665
GSM,H73,5185.6665
GSM,I73,5185.6665
GSM,J73,0
GSM,K73,5185.6665
GSM,L73,0
GSM,M73,-4249.6665
GSM,H94,31161
GSM,I94,0
GSM,J94,31161
GSM,K94,0
GSM,L94,31161
GSM,N94,-32513
LTW,H99,0
LTW,I99,0
LTW,J99,0
LTW,K99,0
LTW,L99,0
LTW,M99,
LTW,H112,0
LTW,I112,0
LTW,J112,0
LTW,K112,0
LTW,L112,
LTW,N112,0
LTW,H132,8953.777
LTW,I132,11214
LTW,J132,4433.3335
LTW,K132,33642
LTW,L132,13300.0005
LTW,M132,-31112
LTW,N132,-18360.0005
LTW,H136,3913.5
LTW,I136,2374.75
LTW,J136,6991
LTW,K136,7124.25
LTW,L136,20973
LTW,M136,-4574.25
LTW,N136,-26073
WVR,H41,19735.666
WVR,I41,19738.5
WVR,J41,19734.25
WVR,K41,19738.5
WVR,L41,19734.25
WVR,M41,-19348.5
WVR,N41,-20514.25
WVR,H55,5226
WVR,I55,4849.3335
WVR,J55,6356
WVR,K55,4849.3335
WVR,L55,6356
WVR,M55,-4572.3335
WVR,N55,-6910
WVR,H80,49999
WVR,I80,0
WVR,J80,49999
WVR,K80,0
WVR,L80,49999
WVR,N80,-51995
WVR,H83,49999
WVR,I83,0
WVR,J83,49999
WVR,K83,0
WVR,L83,149997
WVR,N83,-159981
"@

$rows = $updates -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$header = $true
$ws = $null
$currentSheet = $null

foreach ($line in $rows) {
    if ($header) { $header = $false; continue }  # skip CSV header row
    $parts = $line.Split(",")
    $sheetName = $parts[0]
    $cellRef   = $parts[1]
    $value     = $parts[2]

    if ($sheetName -ne $currentSheet) {
        $ws = $wb.Sheets.Item($sheetName)
        $currentSheet = $sheetName
    }

    if ([string]::IsNullOrEmpty($value)) {
        # Cell dropped from this pull (no market data) -- clear it so it
        # disappears from the sheet rather than showing a stale 0/blank.
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = [double]$value
    }
}
